$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "67.557.31"
$ws.Range("E2").Value = "  +0.05%  "

Set-TextValue "D3" "2.629.75"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue "D5" "595.20"
$ws.Range("E5").Value = "  -0.44%  "

Set-TextValue "D6" "168.57"
$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -2.05%  "

Set-TextValue "D9" "2.629.03"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("E10").Value = "  -1.82%  "

$ws.Range("E12").Value = "  +1.91%  "

$ws.Range("E13").Value = "  +0.00%  "

Set-TextValue "D14" "27.70"
$ws.Range("E14").Value = "  -0.35%  "

Set-TextValue "D15" "3.110.03"
$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("E16").Value = "  -1.05%  "

Set-TextValue "D17" "67.200.67"
$ws.Range("E17").Value = "  -0.34%  "

Set-TextValue "D18" "2.633.30"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("E19").Value = "  +2.63%  "

$ws.Range("E20").Value = "  +4.60%  "

Set-TextValue "D21" "357.08"
$ws.Range("E21").Value = "  -1.77%  "

Set-TextValue "D22" "4.33"

$ws.Range("E24").Value = "  -3.99%  "

$ws.Range("E25").Value = "  -0.04%  "

Set-TextValue "D26" "10.33"
$ws.Range("E26").Value = "  +3.46%  "

Set-TextValue "D27" "69.69"
$ws.Range("E27").Value = "  -1.62%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  -1.35%  "

Set-TextValue "D31" "548.20"
$ws.Range("E31").Value = "  -1.72%  "

Set-TextValue "D32" "7.93"
$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("E33").Value = "  -2.67%  "

Set-TextValue "D34" "1.90"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("E35").Value = "  +4.64%  "

$ws.Range("E36").Value = "  +0.05%  "

Set-TextValue "D38" "156.64"
$ws.Range("E38").Value = "  +1.19%  "

Set-TextValue "D39" "19.04"
$ws.Range("E39").Value = "  -2.52%  "

$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D42" "18.29"
$ws.Range("E42").Value = "  +1.83%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D43" "5.23"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("E46").Value = "  -0.11%  "

Set-TextValue "D47" "153.01"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("E48").Value = "  -1.76%  "

$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("E51").Value = "  -1.12%  "

